# everyday_ver2.xlsx daily update: 2021/12/02 - everyday data updated
# For every sheet, a brand-new row is inserted at row 2 (pushing all existing
# rows down by one, carrying their formatting with them, exactly like a
# manual Excel "Insert Row" above the prior top data row), then the new
# row 2 is filled in with that day's figures.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 台指期換倉成本計算 (columns A-F) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows("2:2").Insert()
$ws1.Range("A2").Value = "日期：2021/12/02"
# B column holds a contract-month code ("202201") that must stay text, not
# be auto-coerced to a number by Excel's input parser.
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "202201"
$ws1.Range("C2").Value = 17649
$ws1.Range("D2").Value = 8583
$ws1.Range("E2").Value = 19484496
$ws1.Range("F2").Value = 17632

# --- Sheet 2: 散戶多空力道 (columns A-B) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows("2:2").Insert()
$ws2.Range("A2").Value = "日期：2021/12/02"
$ws2.Range("B2").Value = 0.01

# --- Sheet 3: 三大法人買賣金額 (columns A-C) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows("2:2").Insert()
$ws3.Range("A2").Value = "110年12月02日"
$ws3.Range("B2").Value = 126.33
$ws3.Range("C2").Value = 12.56

# --- Sheet 4: 大盤多空點位 (columns A-B) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows("2:2").Insert()
$ws4.Range("A2").Value = "110年12月02日"
$ws4.Range("B2").Value = 17670.86

# --- Sheet 5: 期貨大額交易人未沖銷部位 (columns A-N) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows("2:2").Insert()
# A column holds a plain "yyyy/mm/dd" string; without forcing text format
# Excel's input parser would silently reinterpret it as a date serial.
$ws5.Range("A2").NumberFormat = "@"
$ws5.Range("A2").Value = "2021/12/02"
$ws5.Range("B2").Value = 48787
$ws5.Range("C2").Value = 52680
$ws5.Range("D2").Value = -624
$ws5.Range("E2").Value = -2013
$ws5.Range("F2").Value = 25754
$ws5.Range("G2").Value = 47612
$ws5.Range("H2").Value = -1313
$ws5.Range("I2").Value = -1029
$ws5.Range("J2").Value = -21858
$ws5.Range("K2").Value = -284
$ws5.Range("L2").Value = 689
$ws5.Range("M2").Value = -984
$ws5.Range("N2").Value = 1673
